$wb = $excel.ActiveWorkbook

# Helper: write a literal text value into a cell, making sure that
# strings which "look like" a number (e.g. "0.68", "-8.95", "1.0")
# are stored as text (shared string) rather than being auto-coerced
# into a numeric cell by Excel's normal parsing of Range.Value.
function Set-TextValue {
    param($range, [string]$text)

    if ($text -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        # Numeric-looking text: force a Text format so Excel keeps the
        # literal string, then drop back to the Normal style so no
        # lingering number-format is left on the cell.
        $range.NumberFormat = "@"
        $range.Value = $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

# --- Restricciones_del_follower ---------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $ws.Range("A2") "8.95 - y"
Set-TextValue $ws.Range("B2") "-8.95"
Set-TextValue $ws.Range("C2") "J_0_L0_v"
Set-TextValue $ws.Range("D2") "0.68"
Set-TextValue $ws.Range("E2") "0"
Set-TextValue $ws.Range("F2") "0"

Set-TextValue $ws.Range("A3") "-1.9499999999999993 - x + y"
Set-TextValue $ws.Range("B3") "-1.0500000000000007"
Set-TextValue $ws.Range("C3") "J_0_L0_v"
Set-TextValue $ws.Range("D3") "0.24"
Set-TextValue $ws.Range("E3") "0"
Set-TextValue $ws.Range("F3") "0"

Set-TextValue $ws.Range("A4") "-24.9 + x + 2y"
Set-TextValue $ws.Range("B4") "12.899999999999999"
Set-TextValue $ws.Range("C4") "J_0_LP_v"
Set-TextValue $ws.Range("D4") "0.44"
Set-TextValue $ws.Range("E4") "0"
Set-TextValue $ws.Range("F4") "0"

Set-TextValue $ws.Range("A5") "-19.130000000000003 + 4x - y"
Set-TextValue $ws.Range("B5") "7.050000000000001"
Set-TextValue $ws.Range("C5") "J_Ne_L0_v"
Set-TextValue $ws.Range("D5") "1.0"
Set-TextValue $ws.Range("E5") "0"
Set-TextValue $ws.Range("F5") "0"

# --- Punto_modificado ---------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")

Set-TextValue $ws.Range("A1") "x"
Set-TextValue $ws.Range("B1") "y"
Set-TextValue $ws.Range("A2") "7.0"
Set-TextValue $ws.Range("B2") "8.95"

# --- Vector_bf ------------------------------------------------------------
# NB: sheet lookup by name is case-insensitive in this host, and
# "Vector_bf" / "Vector_BF" differ only by case, so they must be
# addressed by their (1-based) tab position instead of by name.
$ws = $wb.Worksheets.Item(5)

Set-TextValue $ws.Range("A1") "vec_bf"
Set-TextValue $ws.Range("A2") "-0.43999999999999995"

# --- Vector_BF ------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)

Set-TextValue $ws.Range("A1") "vec_BF"
Set-TextValue $ws.Range("A2") "1.0"
Set-TextValue $ws.Range("A3") "3.0"
